# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''43.118.09'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +0.76%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.289.96'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +1.57%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.07%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''251.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.61%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.641'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.87%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''73.75'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +4.32%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.03%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.645'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -0.72%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''39.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -5.27%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.0980'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +2.07%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''59.07'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -0.72%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''7.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +1.54%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = '''  +1.49%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''2.636.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +1.74%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''15.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +3.29%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''0.875'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -1.54%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''2.289.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +1.45%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''42.983.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +0.60%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = '''  +2.40%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = '''  +0.79%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''72.65'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -0.60%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''237.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +0.80%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +4.36%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''  -2.35%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''11.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.89%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  -0.25%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = '''  -1.23%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  -0.94%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  -3.12%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''167.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -0.30%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''21.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +0.38%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  +4.52%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  -1.73%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  +4.94%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''30.91'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +10.24%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  +1.58%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''4.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +9.81%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''4.77'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +1.60%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = '''  -2.48%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''14.23'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +14.39%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  +2.44%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '''  +1.95%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '''  +6.11%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''9.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +3.15%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''61.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -3.95%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''4.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -2.29%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '''  +1.10%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  +0.15%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  -1.98%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''100.14'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +5.76%  '
$ws.Range('E51').Style = 'Normal'
